# Update the "Training Dashboard" sheet: decrement PERIOD TO EXPIRE (col H)
# by 1 day and roll LAST UPDATE (col I) forward from 03-Nov-2025 to 04-Nov-2025
# for each data row (rows 3 through 38).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

for ($row = 3; $row -le 38; $row++) {
    $hCell = $ws.Cells.Item($row, 8)   # column H
    $iCell = $ws.Cells.Item($row, 9)   # column I

    $hCell.Value2 = ([double]$hCell.Value2) - 1
    # Leading apostrophe forces Excel to keep this as literal text instead of
    # auto-converting the date-like string into a date serial number.
    $iCell.Value = "'04-Nov-2025"
}
